# Updates the Price (D) and Volume(1h) (E) columns of the cryptos list.
#
# For D-column values that look numeric (e.g. "1.00", "522.16"), a leading
# apostrophe forces Excel to store them as literal text (matching the
# original inline-string cells) instead of coercing them to numbers.
# ClearFormats() is called afterwards so the quote-prefix doesn't leave a
# stray cell style behind (the source cells carry no explicit style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.985.16"
$ws.Range("E2").Value = "  +1.58%  "

$ws.Range("D3").Value = "2.588.55"
$ws.Range("E3").Value = "  +0.07%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'522.16"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.69%  "

$ws.Range("D6").Value = "'139.68"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.76%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").Value = "'0.566"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.18%  "

$ws.Range("D9").Value = "2.600.46"
$ws.Range("E9").Value = "  -0.16%  "

$ws.Range("D10").Value = "'6.54"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.16%  "

$ws.Range("E11").Value = "  -0.28%  "

$ws.Range("D12").Value = "'0.330"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.20%  "

$ws.Range("E13").Value = "  +3.23%  "

$ws.Range("D14").Value = "3.047.18"
$ws.Range("E14").Value = "  +0.12%  "

$ws.Range("D15").Value = "58.944.16"
$ws.Range("E15").Value = "  +1.56%  "

$ws.Range("D16").Value = "'20.57"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.83%  "

$ws.Range("D17").Value = "2.594.76"
$ws.Range("E17").Value = "  -0.46%  "

$ws.Range("E18").Value = "  -0.80%  "

$ws.Range("D19").Value = "'338.95"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.44%  "

$ws.Range("E20").Value = "  -0.15%  "

$ws.Range("D21").Value = "'10.10"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.09%  "

$ws.Range("D22").Value = "'6.47"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.74%  "

$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("D24").Value = "'66.16"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("E25").Value = "  +1.49%  "

$ws.Range("D26").Value = "'0.404"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.50%  "

$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.55%  "

$ws.Range("E28").Value = "  +0.47%  "

$ws.Range("E29").Value = "  +0.06%  "

$ws.Range("D30").Value = "0.0₃0726"
$ws.Range("E30").Value = "  -3.00%  "

$ws.Range("D31").Value = "'5.90"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -5.98%  "

$ws.Range("E32").Value = "  -0.07%  "

$ws.Range("D33").Value = "'18.71"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.41%  "

$ws.Range("D34").Value = "'149.43"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.07%  "

$ws.Range("E35").Value = "  -1.46%  "

$ws.Range("E36").Value = "  -2.25%  "

$ws.Range("D37").Value = "'36.79"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.41%  "

$ws.Range("E38").Value = "  +0.93%  "

$ws.Range("D39").Value = "'0.828"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.31%  "

$ws.Range("D40").Value = "'0.820"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -6.98%  "

$ws.Range("D41").Value = "'3.53"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.81%  "

$ws.Range("D42").Value = "'0.998"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.11%  "

$ws.Range("D43").Value = "'272.74"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("E44").Value = "  +0.94%  "

$ws.Range("D45").Value = "'0.591"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.09%  "

$ws.Range("D46").Value = "'0.0954"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.15%  "

$ws.Range("D47").Value = "'0.0517"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.37%  "

$ws.Range("D48").Value = "'18.42"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.24%  "

$ws.Range("D49").Value = "1.970.80"
$ws.Range("E49").Value = "  -0.27%  "

$ws.Range("D50").Value = "'0.0221"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.16%  "

$ws.Range("D51").Value = "'4.52"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.17%  "
